# 2021-08-03 - Retirada de comentários
# Removes the company logo picture from the "Lista" sheet and genericises
# the "Nº DEVEMADA" label to just "Nº", then leaves the selection where the
# user left it (C6) when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista")

# Remove the logo image (picture shape) that used to sit over A1:B2.
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}

# Generic label instead of the old company-specific "Nº DEVEMADA".
$ws.Range("C5").Value = "Nº"

# Leave the cursor on C6, matching the saved selection.
$ws.Range("C6").Select()
